$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data table (rows 2-7, columns A-F) with the new values.
$ws.Range("A2").Value = 200
$ws.Range("B2").Value = 12
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 10
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0

$ws.Range("A3").Value = 205
$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 12
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1

$ws.Range("A4").Value = 210
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 12
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0

$ws.Range("A5").Value = 215
$ws.Range("B5").Value = 15
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 14
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1

$ws.Range("A6").Value = 220
$ws.Range("B6").Value = 15
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 13
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0

$ws.Range("A7").Value = 225
$ws.Range("B7").Value = 15
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 11
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
